# Bugfixed evaluation and simulated rt_data for components
#
# A new leading observation (date_of_forecast serial 39400 / y_0 2007) is
# prepended to the forecast table, which pushes every existing row down by
# one; at the same time the simulated y_0_forecast / y_1_forecast values
# (columns C and E) are refreshed for the whole series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is brand new (the table used to end at row 18). Prime its formatting
# by copying the format of the current last row (18) into it before any
# values are written, so the date cell (column A) picks up the same style as
# every other date cell instead of Excel minting a fresh one.
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Full refreshed dataset (date serial, y_0, y_0_forecast, y_1, y_1_forecast)
# for rows 2..19 after the insert/shift.
$data = @(
    @(39400, 2007, 3.145939949069287,   2008, 2.455972543253826),
    @(39765, 2008, 1.769627576887389,   2009, 1.087227286828241),
    @(40130, 2009, -4.774178217057779,  2010, -0.8523446516643385),
    @(40494, 2010, 1.97975191822708,    2011, 1.39052144387346),
    @(40862, 2011, 3.452886745653183,   2012, 2.798447799311043),
    @(41228, 2012, 1.239479831392853,   2013, 2.047133666472267),
    @(41592, 2013, 0.2379616621361214,  2014, 1.062273708599726),
    @(41957, 2014, 1.51977456621637,    2015, 0.6176326357196116),
    @(42321, 2015, 1.470039379455756,   2016, 1.577608035818301),
    @(42689, 2016, 1.638797242243251,   2017, 1.369334405341616),
    @(43053, 2017, 2.161565493242668,   2018, 2.486299099038347),
    @(43418, 2018, 2.214251681313772,   2019, 1.362030665126834),
    @(43783, 2019, 0.6066442151010376,  2020, 0.7878236429522678),
    @(44159, 2020, -4.207901339433196,  2021, -0.4850133725290084),
    @(44525, 2021, 1.099928004397532,   2022, 1.102608990163567),
    @(44890, 2022, 2.310042359896225,   2023, 2.086256540666986),
    @(45254, 2023, 0.0464415346324687,  2024, 0.3491198177708599),
    @(45618, 2024, -0.3101476031197148, 2025, 0.2126457877301924)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
